$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: for the zh-cn and de-de sheets, the rows that
# were "Ready for handoff" now reflect a completed handback that is in sync
# with en-US. Each row gets its "Latest Target File" (F) and "Latest
# Handback File" (G) columns populated - mirroring the source file / latest
# handoff file respectively, each as a hyperlink to the very same target -
# the Status column (C) is updated, and the Latest Handback DateTime (H) is
# stamped.
# ---------------------------------------------------------------------------

function Update-HandbackSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)
    Write-Host "Updating sheet:" $SheetName "at" $HandbackDateTime

    # Snapshot the existing hyperlink addresses keyed by the cell they're
    # anchored to (A2/D2/A3/D3) before we start adding new ones.
    $addrByCell = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addrByCell[$hl.Range.Address()] = $hl.Address
    }

    # Status (column C) -> "Handed back: in sync with en-US"
    $ws.Range("C2:C3").Value = "Handed back: in sync with en-US"

    # Column F (Latest Target File) mirrors column A (Source File Name);
    # Column G (Latest Handback File) mirrors column D (Latest Handoff File).
    $ws.Range("F2").Value = $ws.Range("A2").Text
    $ws.Range("G2").Value = $ws.Range("D2").Text
    $ws.Range("F3").Value = $ws.Range("A3").Text
    $ws.Range("G3").Value = $ws.Range("D3").Text

    # Latest Handback DateTime (column H)
    $ws.Range("H2:H3").Value = $HandbackDateTime
    $ws.Range("H2:H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    # Hyperlink the new F/G cells to the same targets as the A/D hyperlinks
    # in the same row (the handback file is, byte-for-byte, the same file
    # that was handed off - hence "in sync").
    $ws.Hyperlinks.Add($ws.Range("F2"), $addrByCell["`$A`$2"], "", "", $ws.Range("A2").Text) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $addrByCell["`$D`$2"], "", "", $ws.Range("D2").Text) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $addrByCell["`$A`$3"], "", "", $ws.Range("A3").Text) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $addrByCell["`$D`$3"], "", "", $ws.Range("D3").Text) | Out-Null
}

Update-HandbackSheet "zh-cn" "2016-03-25 08:52:35"
Update-HandbackSheet "de-de" "2016-03-25 08:52:49"

# The "Status" text is a shared string reused by the Overview roll-up sheet
# (columns B/C hold the zh-cn/de-de status for each file). Updating the
# shared string's text updates every cell that points at it, so mirror the
# same status text there too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2:C3").Value = "Handed back: in sync with en-US"
